$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose updated value is a plain number-looking string (e.g. "0.998").
# Force a text number format first so Excel stores them as text, matching
# the original inlineStr cells, instead of silently converting to numbers.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply updated cell values
$ws.Range('D2').Value = '62.324.50'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '2.451.13'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('D5').Value = '578.12'
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').Value = '143.36'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '2.446.68'
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').Value = '0.345'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').Value = '26.27'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '2.896.23'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = '62.197.82'
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').Value = '2.447.52'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '10.86'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').Value = '327.74'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').Value = '1.95'
$ws.Range('E23').Value = '  -7.06%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = '65.61'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').Value = '9.30'
$ws.Range('E26').Value = '  +3.18%  '
$ws.Range('D27').Value = '586.43'
$ws.Range('E27').Value = '  -5.45%  '
$ws.Range('D28').Value = '2.574.35'
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('D29').Value = '0.0₃0955'
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  -3.66%  '
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').Value = '4.91'
$ws.Range('E35').Value = '  -3.74%  '
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('E37').Value = '  -3.16%  '
$ws.Range('D38').Value = '0.378'
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = '153.30'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('D40').Value = '5.31'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').Value = '43.05'
$ws.Range('E42').Value = '  +2.02%  '
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '2.44'
$ws.Range('E45').Value = '  -4.57%  '
$ws.Range('D46').Value = '3.64'
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('D47').Value = '141.68'
$ws.Range('E47').Value = '  -3.12%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.606'
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0248'
$ws.Range('E49').Value = '  +11.09%  '
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').Value = '19.71'
$ws.Range('E51').Value = '  -4.57%  '
